$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 174.92857
$ws.Range("I11").Value = 174.92857
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 174.92857
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -34.92857000000001
$ws.Range("H62").Value = 2007.5555
$ws.Range("I62").Value = 1878
$ws.Range("J62").Value = 2266.6667
$ws.Range("K62").Value = 1878
$ws.Range("L62").Value = 2266.6667
$ws.Range("M62").Value = -1254
$ws.Range("N62").Value = -3514.6667
$ws.Range("H65").Value = 2007.5555
$ws.Range("I65").Value = 1878
$ws.Range("J65").Value = 2266.6667
$ws.Range("K65").Value = 9390
$ws.Range("L65").Value = 11333.3335
$ws.Range("M65").Value = -6270
$ws.Range("N65").Value = -17573.3335
$ws.Range("H113").Value = 7469.1177
$ws.Range("I113").Value = 5980.8
$ws.Range("J113").Value = 8089.25
$ws.Range("K113").Value = 5980.8
$ws.Range("L113").Value = 8089.25
$ws.Range("M113").Value = -2726.8
$ws.Range("N113").Value = -14597.25
$ws.Range("H137").Value = 1605.8182
$ws.Range("I137").Value = 1024.5883
$ws.Range("J137").Value = 1971.7778
$ws.Range("K137").Value = 3073.7649
$ws.Range("L137").Value = 5915.3334
$ws.Range("M137").Value = -523.7648999999997
$ws.Range("N137").Value = -11015.3334

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 115
$ws.Range("H5").Value = 233.66667
$ws.Range("I5").Value = 100.5
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 100.5
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 11.5
$ws.Range("N5").Value = -724
$ws.Range("H32").Value = 4923.661
$ws.Range("I32").Value = 4257.1777
$ws.Range("J32").Value = 6687.8823
$ws.Range("K32").Value = 4257.1777
$ws.Range("L32").Value = 6687.8823
$ws.Range("M32").Value = -3970.1777
$ws.Range("N32").Value = -7261.8823
$ws.Range("H61").Value = 9724
$ws.Range("I61").Value = 3275.7144
$ws.Range("J61").Value = 100000
$ws.Range("K61").Value = 3275.7144
$ws.Range("L61").Value = 100000
$ws.Range("M61").Value = -3063.7144
$ws.Range("N61").Value = -100424
$ws.Range("H88").Value = 333335330
$ws.Range("I88").Value = 2996
$ws.Range("J88").Value = 500001500
$ws.Range("K88").Value = 2996
$ws.Range("L88").Value = 500001500
$ws.Range("M88").Value = -2590
$ws.Range("N88").Value = -500002312
$ws.Range("H91").Value = 333335330
$ws.Range("I91").Value = 2996
$ws.Range("J91").Value = 500001500
$ws.Range("K91").Value = 2996
$ws.Range("L91").Value = 500001500
$ws.Range("M91").Value = -1592
$ws.Range("N91").Value = -500004308
$ws.Range("H122").Value = 1427986.8
$ws.Range("I122").Value = 1605735.1
$ws.Range("J122").Value = 5999.5
$ws.Range("K122").Value = 4817205.300000001
$ws.Range("L122").Value = 17998.5
$ws.Range("M122").Value = -4814755.300000001
$ws.Range("N122").Value = -22898.5
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6530
$ws.Range("H136").Value = 9724
$ws.Range("I136").Value = 3275.7144
$ws.Range("J136").Value = 100000
$ws.Range("K136").Value = 9827.143199999999
$ws.Range("L136").Value = 300000
$ws.Range("M136").Value = -7277.143199999999
$ws.Range("N136").Value = -305100

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 233.66667
$ws.Range("I4").Value = 100.5
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 100.5
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = 14.5
$ws.Range("N4").Value = -730

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 388.9
$ws.Range("I22").Value = 270
$ws.Range("J22").Value = 666.3333
$ws.Range("K22").Value = 270
$ws.Range("L22").Value = 666.3333
$ws.Range("M22").Value = 80
$ws.Range("N22").Value = -1366.3333
$ws.Range("H99").Value = 17882858
$ws.Range("I99").Value = 30000
$ws.Range("J99").Value = 31272500
$ws.Range("K99").Value = 30000
$ws.Range("L99").Value = 31272500
$ws.Range("M99").Value = -28502
$ws.Range("N99").Value = -31275496
$ws.Range("H118").Value = 39999.848
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 39999.848
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 39999.848
$ws.Range("N118").Value = -43313.848
$ws.Range("H122").Value = 2596
$ws.Range("I122").Value = 2245
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 6735
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -4285
$ws.Range("N122").Value = -16900
$ws.Range("H126").Value = 17882858
$ws.Range("I126").Value = 30000
$ws.Range("J126").Value = 31272500
$ws.Range("K126").Value = 90000
$ws.Range("L126").Value = 93817500
$ws.Range("M126").Value = -87530
$ws.Range("N126").Value = -93822440
$ws.Range("H132").Value = 2531.5806
$ws.Range("I132").Value = 1658.12
$ws.Range("J132").Value = 6171
$ws.Range("K132").Value = 4974.36
$ws.Range("L132").Value = 18513
$ws.Range("M132").Value = -2444.36
$ws.Range("N132").Value = -23573
$ws.Range("H134").Value = 2637.7097
$ws.Range("I134").Value = 2865.5186
$ws.Range("J134").Value = 1100
$ws.Range("K134").Value = 8596.5558
$ws.Range("L134").Value = 3300
$ws.Range("M134").Value = -6061.5558
$ws.Range("N134").Value = -8370

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 102.05882
$ws.Range("I14").Value = 102.05882
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 306.17646
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -133.17646
$ws.Range("H68").Value = 2953.2097
$ws.Range("I68").Value = 4230.067
$ws.Range("J68").Value = 1756.1562
$ws.Range("K68").Value = 12690.201
$ws.Range("L68").Value = 5268.4686
$ws.Range("M68").Value = -11879.201
$ws.Range("N68").Value = -6890.4686
$ws.Range("H71").Value = 2953.2097
$ws.Range("I71").Value = 4230.067
$ws.Range("J71").Value = 1756.1562
$ws.Range("K71").Value = 38070.603
$ws.Range("L71").Value = 15805.4058
$ws.Range("M71").Value = -34014.603
$ws.Range("N71").Value = -23917.4058
$ws.Range("H97").Value = 8333816.5
$ws.Range("I97").Value = 16667033
$ws.Range("J97").Value = 600
$ws.Range("K97").Value = 50001099
$ws.Range("L97").Value = 1800
$ws.Range("M97").Value = -50000603
$ws.Range("N97").Value = -2792
$ws.Range("H113").Value = 208806.56
$ws.Range("I113").Value = 454.14706
$ws.Range("J113").Value = 714805.3
$ws.Range("K113").Value = 1362.44118
$ws.Range("L113").Value = 2144415.9
$ws.Range("M113").Value = 807.55882
$ws.Range("N113").Value = -2148755.9
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H118").Value = 1179.8
$ws.Range("I118").Value = 749.75
$ws.Range("J118").Value = 2900
$ws.Range("K118").Value = 2249.25
$ws.Range("L118").Value = 8700
$ws.Range("M118").Value = -1006.25
$ws.Range("N118").Value = -11186
$ws.Range("H124").Value = 750
$ws.Range("I124").Value = 500
$ws.Range("J124").Value = 1000
$ws.Range("K124").Value = 1500
$ws.Range("L124").Value = 3000
$ws.Range("M124").Value = 3410
$ws.Range("N124").Value = -12820
$ws.Range("H125").Value = 3019.8
$ws.Range("I125").Value = 1200
$ws.Range("J125").Value = 4233
$ws.Range("K125").Value = 3600
$ws.Range("L125").Value = 12699
$ws.Range("M125").Value = 1320
$ws.Range("N125").Value = -22539

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 30000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 30000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180
$ws.Range("H122").Value = 19137468
$ws.Range("I122").Value = 12964962
$ws.Range("J122").Value = 50000000
$ws.Range("K122").Value = 38894886
$ws.Range("L122").Value = 150000000
$ws.Range("M122").Value = -38892436
$ws.Range("N122").Value = -150004900

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7938
$ws.Range("I7").Value = 3034.6667
$ws.Range("J7").Value = 10880
$ws.Range("K7").Value = 3034.6667
$ws.Range("L7").Value = 10880
$ws.Range("M7").Value = -2922.6667
$ws.Range("N7").Value = -11104
$ws.Range("I40").Value = 250002850
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 250002850
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -250002714
$ws.Range("N40").Value = -4272
$ws.Range("H82").Value = 11466143
$ws.Range("I82").Value = 21320
$ws.Range("J82").Value = 21003496
$ws.Range("K82").Value = 21320
$ws.Range("L82").Value = 21003496
$ws.Range("M82").Value = -20959
$ws.Range("N82").Value = -21004218
$ws.Range("H85").Value = 11466143
$ws.Range("I85").Value = 21320
$ws.Range("J85").Value = 21003496
$ws.Range("K85").Value = 21320
$ws.Range("L85").Value = 21003496
$ws.Range("M85").Value = -20072
$ws.Range("N85").Value = -21005992
$ws.Range("H96").Value = 19299
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 19299
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 19299
$ws.Range("N96").Value = -24791
$ws.Range("H126").Value = 7938
$ws.Range("I126").Value = 3034.6667
$ws.Range("J126").Value = 10880
$ws.Range("K126").Value = 9104.000100000001
$ws.Range("L126").Value = 32640
$ws.Range("M126").Value = -6634.000100000001
$ws.Range("N126").Value = -37580
$ws.Range("H136").Value = 6047.2334
$ws.Range("I136").Value = 3000.3809
$ws.Range("J136").Value = 13156.556
$ws.Range("K136").Value = 9001.1427
$ws.Range("L136").Value = 39469.66800000001
$ws.Range("M136").Value = -6451.1427
$ws.Range("N136").Value = -44569.66800000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 22144
$ws.Range("I99").Value = 21432
$ws.Range("J99").Value = 22500
$ws.Range("K99").Value = 21432
$ws.Range("L99").Value = 22500
$ws.Range("M99").Value = -18437
$ws.Range("N99").Value = -28490
$ws.Range("H116").Value = 45000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 45000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 45000
$ws.Range("N116").Value = -54178
$ws.Range("H122").Value = 4785.364
$ws.Range("I122").Value = 4472.3335
$ws.Range("J122").Value = 5161
$ws.Range("K122").Value = 13417.0005
$ws.Range("L122").Value = 15483
$ws.Range("M122").Value = -10967.0005
$ws.Range("N122").Value = -20383
$ws.Range("H126").Value = 1562.625
$ws.Range("I126").Value = 1034
$ws.Range("J126").Value = 1879.8
$ws.Range("K126").Value = 3102
$ws.Range("L126").Value = 5639.4
$ws.Range("M126").Value = -632
$ws.Range("N126").Value = -10579.4

Write-Output "All edits applied"